$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.762.80"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.633.46"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0633"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D13").Value = "1.859.76"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "1.634.62"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.559"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "25.784.05"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E34").Value = "  -1.33%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  +1.83%  "
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.801"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").Value = "1.768.44"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  +3.42%  "
